$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (no auto numeric/date coercion, no style change).
# Stage the text in a scratch cell formatted as Text, copy it, then PasteSpecial
# (values only) into the destination so the destination keeps its original style.
function Set-TextValue($addr, $val) {
    $scratch = $ws.Range("Z1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $val
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue "D2" '29.007.66'
Set-TextValue "E2" '  +1.17%  '
Set-TextValue "D3" '1.886.01'
Set-TextValue "E3" '  +0.59%  '
Set-TextValue "E4" '  +0.07%  '
Set-TextValue "D5" '331.27'
Set-TextValue "E5" '  -2.11%  '
Set-TextValue "D6" '1.001'
Set-TextValue "E6" '  +0.17%  '
Set-TextValue "D7" '0.4591'
Set-TextValue "E7" '  -2.46%  '
Set-TextValue "D8" '0.4081'
Set-TextValue "E8" '  +1.91%  '
Set-TextValue "D9" '47.18'
Set-TextValue "E9" '  -1.25%  '
Set-TextValue "D10" '0.07995'
Set-TextValue "E10" '  -0.76%  '
Set-TextValue "D11" '0.9894'
Set-TextValue "E11" '  -1.66%  '
Set-TextValue "E12" '  -2.22%  '
Set-TextValue "D13" '1.875.23'
Set-TextValue "E13" '  +0.76%  '
Set-TextValue "D14" '5.896'
Set-TextValue "E14" '  -2.88%  '
Set-TextValue "D15" '7.055'
Set-TextValue "E15" '  -3.32%  '
Set-TextValue "D16" '1.001'
Set-TextValue "E16" '  -0.06%  '
Set-TextValue "D17" '88.66'
Set-TextValue "E17" '  -2.22%  '
Set-TextValue "D18" '0.00001026'
Set-TextValue "E18" '  -1.83%  '
Set-TextValue "D19" '0.06552'
Set-TextValue "E19" '  -0.94%  '
Set-TextValue "D21" '1.000'
Set-TextValue "E21" '  +0.13%  '
Set-TextValue "D22" '29.048.23'
Set-TextValue "E22" '  +1.25%  '
Set-TextValue "D23" '5.397'
Set-TextValue "E23" '  -2.13%  '
Set-TextValue "D24" '11.24'
Set-TextValue "E24" '  +1.57%  '
Set-TextValue "D25" '2.208'
Set-TextValue "E25" '  -2.35%  '
Set-TextValue "D26" '2.122.57'
Set-TextValue "E26" '  +1.80%  '
Set-TextValue "D27" '156.93'
Set-TextValue "E27" '  -2.33%  '
Set-TextValue "D28" '19.61'
Set-TextValue "E28" '  -1.17%  '
Set-TextValue "D29" '2.098'
Set-TextValue "E29" '  -1.70%  '
Set-TextValue "D30" '5.413'
Set-TextValue "E30" '  -1.97%  '
Set-TextValue "D31" '117.57'
Set-TextValue "E31" '  -2.15%  '
Set-TextValue "D32" '0.9734'
Set-TextValue "E32" '  -1.41%  '
Set-TextValue "D33" '0.09314'
Set-TextValue "E33" '  -2.50%  '
Set-TextValue "D34" '3.601'
Set-TextValue "E34" '  -1.96%  '
Set-TextValue "D35" '1.404'
Set-TextValue "E35" '  +1.08%  '
Set-TextValue "D36" '5.269'
Set-TextValue "E36" '  -2.17%  '
Set-TextValue "D37" '0.06036'
Set-TextValue "E37" '  -2.70%  '
Set-TextValue "D38" '0.02222'
Set-TextValue "E38" '  -2.06%  '
Set-TextValue "D39" '8.245'
Set-TextValue "E39" '  -2.89%  '
Set-TextValue "D40" '1.182'
Set-TextValue "E40" '  -0.49%  '
Set-TextValue "D41" '1.000'
Set-TextValue "E41" '  +0.12%  '
Set-TextValue "D42" '0.5757'
Set-TextValue "E42" '  -3.60%  '
Set-TextValue "D43" '0.1817'
Set-TextValue "E43" '  -3.95%  '
Set-TextValue "D44" '10.10'
Set-TextValue "E44" '  -2.55%  '
Set-TextValue "E45" '  -0.28%  '
Set-TextValue "D46" '12.00'
Set-TextValue "E46" '  -2.25%  '
Set-TextValue "D47" '2.253'
Set-TextValue "E47" '  +8.28%  '
Set-TextValue "D48" '0.5452'
Set-TextValue "E48" '  -2.47%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D49" '0.07026'
Set-TextValue "E49" '  -5.30%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D50" '1.892'
Set-TextValue "E50" '  -3.79%  '
Set-TextValue "D51" '45.53'
Set-TextValue "E51" '  +13.65%  '

$ws.Range("Z1").Clear()
$excel.CutCopyMode = $false
